# Applies numeric corrections produced by the scheduled profit-recalculation
# runner across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 909.94446
$ws.Range("I6").Value = 698.6
$ws.Range("J6").Value = 1966.6666
$ws.Range("K6").Value = 2095.8
$ws.Range("L6").Value = 5899.9998
$ws.Range("M6").Value = -1983.8
$ws.Range("N6").Value = -6123.9998

# Row 8
$ws.Range("H8").Value = 406.375
$ws.Range("I8").Value = 406.375
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1219.125
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1080.125
$ws.Range("N8").Value = ""

# Row 98
$ws.Range("H98").Value = 2794.4048
$ws.Range("I98").Value = 723.0909
$ws.Range("K98").Value = 723.0909
$ws.Range("M98").Value = 774.9091

# Row 107
$ws.Range("H107").Value = 1126
$ws.Range("I107").Value = 1168
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1168
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 752
$ws.Range("N107").Value = -4840

# Row 116
$ws.Range("H116").Value = 884086.9
$ws.Range("I116").Value = 1351979.6
$ws.Range("K116").Value = 1351979.6
$ws.Range("M116").Value = -1348537.6

# Row 122
$ws.Range("H122").Value = 2794.4048
$ws.Range("I122").Value = 723.0909
$ws.Range("K122").Value = 2169.2727
$ws.Range("M122").Value = 280.7273

# Row 132
$ws.Range("H132").Value = 3191418
$ws.Range("I132").Value = 3682002.5
$ws.Range("K132").Value = 11046007.5
$ws.Range("M132").Value = -11043477.5

# Row 137
$ws.Range("H137").Value = 23881.715
$ws.Range("I137").Value = 30727.818
$ws.Range("J137").Value = 19451.883
$ws.Range("K137").Value = 92183.454
$ws.Range("L137").Value = 58355.649
$ws.Range("M137").Value = -89633.454
$ws.Range("N137").Value = -63455.649

# Row 138
$ws.Range("H138").Value = 43158.27
$ws.Range("I138").Value = 3288.75
$ws.Range("J138").Value = 106949.5
$ws.Range("K138").Value = 9866.25
$ws.Range("L138").Value = 320848.5
$ws.Range("M138").Value = -4726.25
$ws.Range("N138").Value = -331128.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15389.312
$ws.Range("I32").Value = 15389.312
$ws.Range("K32").Value = 15389.312
$ws.Range("M32").Value = -15102.312

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1978
$ws.Range("I86").Value = 1606.25
$ws.Range("J86").Value = 2349.75
$ws.Range("K86").Value = 1606.25
$ws.Range("L86").Value = 2349.75
$ws.Range("M86").Value = -483.25
$ws.Range("N86").Value = -4595.75

# Row 89
$ws.Range("H89").Value = 1978
$ws.Range("I89").Value = 1606.25
$ws.Range("J89").Value = 2349.75
$ws.Range("K89").Value = 8031.25
$ws.Range("L89").Value = 11748.75
$ws.Range("M89").Value = -2415.25
$ws.Range("N89").Value = -22980.75

# Row 105
$ws.Range("H105").Value = 2009.7667
$ws.Range("I105").Value = 1695.96
$ws.Range("J105").Value = 3578.8
$ws.Range("K105").Value = 1695.96
$ws.Range("L105").Value = 3578.8
$ws.Range("M105").Value = 51.03999999999996
$ws.Range("N105").Value = -7072.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1192322.2
$ws.Range("I31").Value = 1819852.9
$ws.Range("K31").Value = 1819852.9
$ws.Range("M31").Value = -1819557.9

# Row 34
$ws.Range("H34").Value = 1192322.2
$ws.Range("I34").Value = 1819852.9
$ws.Range("K34").Value = 1819852.9
$ws.Range("M34").Value = -1819650.9

# Row 58
$ws.Range("H58").Value = 1186.6285
$ws.Range("I58").Value = 967.93335
$ws.Range("J58").Value = 2498.8
$ws.Range("K58").Value = 967.93335
$ws.Range("L58").Value = 2498.8
$ws.Range("M58").Value = -764.93335
$ws.Range("N58").Value = -2904.8

# Row 105
$ws.Range("H105").Value = 1984.2307
$ws.Range("I105").Value = 1388
$ws.Range("J105").Value = 2938.2
$ws.Range("K105").Value = 1388
$ws.Range("L105").Value = 2938.2
$ws.Range("M105").Value = 359
$ws.Range("N105").Value = -6432.2

# Row 136
$ws.Range("H136").Value = 1186.6285
$ws.Range("I136").Value = 967.93335
$ws.Range("J136").Value = 2498.8
$ws.Range("K136").Value = 2903.80005
$ws.Range("L136").Value = 7496.400000000001
$ws.Range("M136").Value = -353.8000499999998
$ws.Range("N136").Value = -12596.4

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 928.44446
$ws.Range("I114").Value = 411.4
$ws.Range("J114").Value = 1574.75
$ws.Range("K114").Value = 1234.2
$ws.Range("L114").Value = 4724.25
$ws.Range("M114").Value = 2019.8
$ws.Range("N114").Value = -11232.25

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1331.4445
$ws.Range("J113").Value = 1430.3334
$ws.Range("L113").Value = 1430.3334
$ws.Range("N113").Value = -5770.3334

# Row 122
$ws.Range("H122").Value = 15154328
$ws.Range("I122").Value = 2668.6538
$ws.Range("K122").Value = 8005.9614
$ws.Range("M122").Value = -5555.9614

# Row 132
$ws.Range("H132").Value = 2339.7058
$ws.Range("I132").Value = 2140.2593
$ws.Range("J132").Value = 3109
$ws.Range("K132").Value = 6420.777900000001
$ws.Range("L132").Value = 9327
$ws.Range("M132").Value = -3890.777900000001
$ws.Range("N132").Value = -14387

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2361.0588
$ws.Range("I7").Value = 2295.5715
$ws.Range("J7").Value = 2666.6667
$ws.Range("K7").Value = 2295.5715
$ws.Range("L7").Value = 2666.6667
$ws.Range("M7").Value = -2183.5715
$ws.Range("N7").Value = -2890.6667

# Row 40
$ws.Range("H40").Value = 4157.4165
$ws.Range("I40").Value = 3988.9
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3988.9
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3852.9
$ws.Range("N40").Value = -5272

# Row 61
$ws.Range("H61").Value = 1537.579
$ws.Range("I61").Value = 1310.4615
$ws.Range("J61").Value = 2029.6666
$ws.Range("K61").Value = 1310.4615
$ws.Range("L61").Value = 2029.6666
$ws.Range("M61").Value = -1108.4615
$ws.Range("N61").Value = -2433.6666

# Row 113
$ws.Range("H113").Value = 1537.579
$ws.Range("I113").Value = 1310.4615
$ws.Range("J113").Value = 2029.6666
$ws.Range("K113").Value = 1310.4615
$ws.Range("L113").Value = 2029.6666
$ws.Range("M113").Value = 859.5385000000001
$ws.Range("N113").Value = -6369.6666

# Row 126
$ws.Range("H126").Value = 2361.0588
$ws.Range("I126").Value = 2295.5715
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 6886.7145
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -4416.7145
$ws.Range("N126").Value = -12940.0001

# Row 132
$ws.Range("H132").Value = 3367.842
$ws.Range("I132").Value = 2838.1538
$ws.Range("K132").Value = 8514.4614
$ws.Range("M132").Value = -5984.4614

# Row 136
$ws.Range("H136").Value = 4304.091
$ws.Range("I136").Value = 4435.294
$ws.Range("J136").Value = 3858
$ws.Range("K136").Value = 13305.882
$ws.Range("L136").Value = 11574
$ws.Range("M136").Value = -10755.882
$ws.Range("N136").Value = -16674

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 38835.8
$ws.Range("I132").Value = 44930.06
$ws.Range("J132").Value = 4301.6665
$ws.Range("K132").Value = 134790.18
$ws.Range("L132").Value = 12904.9995
$ws.Range("M132").Value = -132260.18
$ws.Range("N132").Value = -17964.9995

# Row 136
$ws.Range("H136").Value = 23484.1
$ws.Range("I136").Value = 25101.893
$ws.Range("K136").Value = 75305.679
$ws.Range("M136").Value = -72755.679
Write-Output "Applied profit recalculation updates"
